$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "E2"  = 183.33
    "E4"  = 700
    "E6"  = 800
    "E7"  = 800
    "E8"  = 800
    "E10" = 800
    "E11" = 800
    "E14" = 700
    "E15" = 700
    "E16" = 800
    "E17" = 600
    "E18" = 600
    "E20" = 300
    "E21" = 450
    "E22" = 500
    "E23" = 800
    "E24" = 600
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
